$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$values = @(
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/kraljevica-marka/637d05681adea1f4a70d45ad',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/petrovaradin/637ceb9e0dc92027c0053ff6',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/podbara/63936deb3dcdd44a11062fc5',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/zeleznicka-stanica/63a4954914d7bb44f10bc78f',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/centar-novi-sad/639c8aa3647bf0c5ef099c2d',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/nova-detelinara-novi-sad/63a5fe9dd6d88ed09b0e9590',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/detelinara/63947a38e0c0f5cd0302d49a',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/nova-23/63a1ef80ecf93efbdf032075',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novi-sad-centar/63a0687a4ba299f4690c899a',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/liman-4/638fa1332e1f67b5cc0d4553',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/grbavica/639c57715032de31920decf6',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/grbavica/638e29caa7ac6955de09b0f4',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-evrope/6394a6fedda11d36da0d2c08',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/liman-4/63a30adf6ad43758790d67f3',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/adice/638f44e561a0886cf8066425',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novi-sad-centar/63a5933c976a14e1b5072415',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/socijalno-novi-sad/63a1b5bbc8c628ea2e07d129',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/telep/63a2fe126d802c181007b8c5',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/alekse-santica/639bfc905da2ca55470cd3aa',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/presernova/639751c8fc647bd628054cb6',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/kraljevica-marka-1/63a6ff6af9b6f0b02d05f861',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-novi-sad/63a35b9b998876901a032154',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novi-sad-centar/63a44019f51974957103ae1b',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novo-naselje/636a4ce748169f3f6b0c6b17',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/telep/636ced22af453926ea086053',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/cara-dusana/637e37819a0b5b0651069d8b',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novi-sad-centar/638dd54c107c506f5e039a28',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novi-sad-centar/63a6e38c8e40023ab400aeb0',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-patrijarha-pavla-3/5fc7cb6f0c7cde55974ee3d7',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/adice-novi-sad/63a458ded20f9e794102a8bb',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/jevrejska-ulica/63356c9d9b41e5ba5b066ff0',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/veternik/639c72c97d2296103909291e',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/petrovaradin/63972c5f7c19077c1d0aa6a1',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/karadjordjeva/63246aaf8e1d4f08f80b7aba',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/dragise-brasovana-10/63a59f034ab49938530bb75f',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/sime-matavulja/639a052cd7a2ed61540befdb',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/socijalno/63a44f9a14a25d684d063a2d',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-oslobodjenja/6389dd62c4229d01930b9237',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-evrope/63a2eff9934fd2a86b049bd4',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/rotkvarija/5edfe04a0c7cde06361a6862',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/podbara/6397399cd4250c7b4702402a',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/janka-cmelika/63a841b4d98fd06f7401ca50',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/narodnog-fronta/639ca7e4b30b2e1e0700c120',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/nova-detelinara/6389fbf230a29508bf00d0b2',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/sajmiste/639da74842fba669a4021171',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/liman-2/63a200f5cc4142344a078172',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/nova-detelinara/638a452d3745bb956509b53a',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-oslobodjenja/639dc2eacb93fbf59f0eb9ee',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/socijalno/6358f13a3aef3aead509a448',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/temerinska/63a1ef7becf93efbdf032069',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/desanke-maksimovic-56/63a42154e23e98ba580f9815',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/mise-dimitrijevica/63a19bf70a2cf277ed09a698',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novi-sad-centar/638f5a6502f74e104a029420',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/centar-novi-sad/6380f745998f51c35008f533',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/podbara/63a45be4d20f9e794102a9ab',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novo-naselje-novi-sad/636a2b48e99ebe5ff60baab1',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/sajam-novi-sad/637cd889a447fc5d61016a4f',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-evrope/63a582433aebf5efca0aeb1b',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-evrope/61f924e8c6588b0e4c6a71ed',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/petrovaradin-petrovaradin/6308b89344c2131d56095849',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/vojvode-supljikca/638a32cc8e4d389d66092f92',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/djordja-servickog/638212df70883b347508f73b',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/petefi-sandora/63a1cb0710ecf1b16f099627',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-oslobodjenja/62c81eaaf9d0ff6fd4146d63',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/brace-ribnikar/62a0abef0480e5590d4e8383',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/centar-novi-sad/6389dd8348993d1ebc0fc36e',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/telep/63a4402ff51974957103ae26',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/studenicka/638dd420107c506f5e0399b7',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-evrope/6394867be022754eca05892a',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/riblja-pijaca/62e6ac6dba88fa6bf9024b0c',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/podbara/6398b5fc84f134d14603f80f',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/puskinova/63a5e4020cd7e44acf05b52b',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/hadzi-ruvimova/639df51e4fc9ba36190b9a80',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-oslobodjenja/63a5cabb187d7566850f01e2',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/futoska/639c39aee1a14cb9f90bb0c3',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/cara-dusana/6399b38f063fdcc3990a1d57',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-oslobodjenja/63a18fc4d1ca73d2c3017a0d',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-despota-stefana/6141fc53aa67a071b6193833',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novosadskog-sajma/61cb06fedf14b705d25be49a',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/petrovaradin/637d850acfd3a0af160c234d',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/laze-pacua/637b98a6c60c1f1ff5030307',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/cara-dusana/6374f309157e4973e4026f93',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/satelit/63a0c8dd53f1afab8d08f2fd',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/petrovaradin/636923bbf202e19cec0ebc8d',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-evrope/63a5c64a5231129cd70e9489',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/centar-novi-sad/63a1986dd1ca73d2c3017a5f',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/janka-cmelika/6384993710b7c9a9fa033d54',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/kraljevica-marka/637de4d9dd2c11f56f01d589',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-evrope-novi-sad/63a33f7c2a81f487d902f803',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/grbavica/63a02dcc7717fc2ac50e2551',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/sime-matavulja/63a5e3ec0cd7e44acf05b50d',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/toplice-milana/63a064e48bf6e0cd67013250',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novi-sad-centar/639c71a17d22961039092880',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/djordja-niksica-johana/62bc53fd01301b51b60fc5e4',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-oslobodjenja/63a487c6e7a5758c420a18b5',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/jovana-hranilovica/63a5d9a263fe7fc2ba056372',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/starine-novaka/639dde307d135708930afde5',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novi-sad-centar/63a334df2a81f487d902f770',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/kosovska/639a1aa7cf4cd3eedd0bea2d',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/telep/62e5094a9061321ebd0fb78d',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/centar-novi-sad/6393887ee035c2e1620176ab',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-evrope/639c55eab5308ad1700fe069',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/rotkvarija/637e1fabb9413c60070a0a63',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/grbavica/6389a7eaf1ac88e06e09d8c0',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/slovacka/63a492a514d7bb44f10bc6d9',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/socijalno/639c55e75032de31920decc7',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/olge-hadzic-14/6090611360efac3b9f64feea',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/zeleznicka-stanica/63a4256359f2a61b410cfcad',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/sajam-novi-sad/639b7165daa71d31090b8bd1',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/liman-3/63989b9abc74447ddd005b9d',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/liman-2/63873d629d960dbaa10f1855',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/pasterova/63a03330b98694ea2d07bfe8',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novi-sad-centar/638207da70883b347508f68e',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/branka-radicevica/63a2fb446d802c181007b807',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-evrope/63a325ce4da8c17be00a6cd1',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/brace-lucic/633c25466422d9a1b909c8c5',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/radnicka/63a02ce07717fc2ac50e2504',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novi-sad-centar/639c71a37d22961039092886',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/sajmiste/639dfb4bcdc2aba3fe066fa8',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/draska-redjepa/638dd9a7f6d96eda4009cb98',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/nova-detelinara/638a1aff47612c11f40ef5e4',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/podbara/63a2288f89aad9162907c6a5',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/socijalno-novi-sad/63863d6027c35623870a42c1',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/petra-drapsina/6385e2a8014c85dda70808f4',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/andje-rankovic/62bb4049f0c35f4daa2c5f03',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novo-naselje/6390cac8f6ebce13f0055e7b',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/avijacija-novi-sad/638f236d09cd16398f066dc5',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novi-sad-centar/638f1ab1d1cca64dd409f78d',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/adice-novi-sad/639760c65c4f3faf260e0765',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novi-sad-centar/638dc6c25696880d4f035e9c',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/detelinara/63a59120976a14e1b5072372',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/klisa-novi-sad/63878eb965467da6f90331b8',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/zlatariceva-18/5d4b18cee466e627cf5adad8',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/stefana-stefanovica/6389e09148993d1ebc0fc44d',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/podbara/63a80e151541e67af20ce6e8',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/rodina/63a76df84d2b1035a505f5a5',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/strazilovska/63a1a9f1bdc1be63350b8a9a',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/nova-detelinara-novi-sad/63a717d6d0eaf9aaf80f24e9',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/salajka/639b035d9745dd65f009e5e8',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-oslobodjenja/63986362d3b70b22bf001a45',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-oslobodjenja/6367ec09b5b958415c058a01',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/socijalno-novi-sad/636a8dc291f9a3ea0b081f68',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/jovana-hranilovica/63a5e3fc0cd7e44acf05b51f',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/sajmiste/636a3cad03bc008f70061c45',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/slovacka/63a459f1d20f9e794102a90d',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/liman-novi-sad/63a2cefdf9ae2c05e10a38da',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-oslobodjenja/638de24dacc7d52c660f9e4a',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/livade/638a040f5cd44949fb0fffb3',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/nova-detelinara/6399a49dba32d7809b07acda',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/kraljevica-marka/636d403609b597bd140f5eb6',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-jovana-ducica/63a5dc9c4805759e8b03663a',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/adice/639738eed4250c7b47023fed',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/nova-detelinara/639b22749b0294e3c801d52e',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/stevana-doronjskog/63a87988c72a6ec04105b1e1',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/djurdja-brankovica-10/63a329a7c9966cb2d50f437a',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/grbavica-novi-sad/63a5e2800cd7e44acf05b4d0',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/grbavica/63889c51937bf2813c0f1797',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/detelinara/63a0680f4ba299f4690c898c',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/preradoviceva-ulica/6396df856bcaebc19109b153',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/adice/6398716271e4fe0f320f4e90',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/sajmiste/639b674bd934b5892e00675b',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/detelinara/63a325cd4da8c17be00a6ccb',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/kraljevica-marka/63763483fada89156c0818d3',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-cara-lazara/611d12ec5cf47054ac2e6887',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/grbavica/637b61c82b6386ba6309c626',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/olge-hadzic/637ea9dc63856c9beb017d0c',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/jovana-ristica-16/63a8ab7444e8671d9001ee5d',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/adice/63a979d0da780523ce0eb6bf',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-jase-tomica/63a983ba76dfb2d8dc0ab414',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/socijalno/63a9877d76dfb2d8dc0ab50e',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-oslobodjenja/63a95a6ac94c64a20e05757a',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-patrijarha-pavla-somborski-bulevar/63a9cd1c603183211b0751c9',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-patrijarha-pavla-somborski-bulevar/63a9b04651832c9cf50e625c',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/liman-1/63a95b845d5f1a6b8400fc98',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/grbavica-novi-sad/63a9bb4fa9a58aa6100ba735',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/zivojina-culuma/63a997ae6f69f369ea0d7cef',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-oslobodjenja/63a9a19f04bbcdeb3603954f',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-oslobodjenja/63a9bfc1a9a58aa6100ba891',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/gagarinova/63a9678f0887df837a069f8f',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novi-sad-centar/63a9dae6da8414c5bb044583',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/grbavica/62e95cea83c2c35bf400c465',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/mise-dimitrijevica-2/6343d9c75d231820e7018619',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/cara-dusana/63ab2b719d578a4747045772',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/cara-dusana/63aaf16a60f003fd450a16aa',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novi-sad-centar/63ab1547e3f509222d0746d3',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/grbavica/63ab2d96baebe57d4007de06',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/nova-detelinara/63aaca7055777b731807d25b',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-evrope/63aad6ee6ad53dbca60d1805',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/veternicka-rampa/63aac81b49719e791c09dd31',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/brace-ribnikar-21/5efd09150c7cde057d4bd2a2',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/sterijina/59d8bcc27ce3af04392915e9',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/grbavica/63ab02b98a083674c601d260',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/cara-dusana/63aaf4d260f003fd450a17b5',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/molijera/63ab0d8fa438df9f960500ca',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/nova/63aacb4a10bb04ad95045c62',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/branimira-cosica/63a5c018452370e7150dd6b8',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-oslobodjenja/639c4baa17a1189a6d0fc832',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/cirpanova/63aadfc0e95f94c25201c408',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novi-sad-centar/63ab1d7587f2b3deda0c56fc',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/socijalno-novi-sad/63aab80586a6c86e1b0dbb4b',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/adice/63ab3995baebe57d4007df3d',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/svetosavska/63aa2cb7008faef69d04aad8',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/detelinara-novi-sad/63ab0c5e2bf1f7b6b0053615',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/kisacka/63aad63b6ad53dbca60d17be',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/masarikova/63aad6306ad53dbca60d17b6',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/lasla-gala/63aa2cbb008faef69d04aade',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/jugovicevo/63ab0afb8bb85133750c2f9c',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/avijaticarsko-naselje/63aacaf16ad53dbca60d16f5',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/grbavica/63aaf4d460f003fd450a17b8',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novo-naselje/63aad909e082d0813103f715',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/grbavica/63ac4432303617380309956c',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/cirpanova/63ac7b982af8a275fc0fbb5e',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/vuka-karadzica/63accfcf6422a9c9f50b0fad',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/draska-redjepa/637e37819a0b5b0651069d8b',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/grbavica/63ac70942af8a275fc0fbaa8',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-evrope/63a9bb925475f663e303052e',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/molijerova/63ac43ea303617380309954b',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/tozin-sokak/63ab3e6524b3290f0e0db48a',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/sajmiste/63ac8c09b65bd1bae40c0f64',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/kolo-srpskih-sestara/63ac42d81879b258f90aa22a',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/grbavica-novi-sad/63ab44996c5999e26d0720a0',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/cirpanova/63ab54b6d4d115f4e90ba1a9',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-evrope/63ac70972af8a275fc0fbaab',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-evrope-novi-sad/63ac9621b65bd1bae40c0ffd',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/grbavica-novi-sad/63ac5de369d41fedf9045a09',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/liman-2-novi-sad/63ac7a0b2af8a275fc0fbafb',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/nova-detelinara/63ac0e30bae24ac33b05a424',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/nova-detelinara/63ac2a270d6db319d702b574',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/cankareva/63b2ee9fb65028de5608665f',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/liman-4-novi-sad/63aec89f2562e3893a016eb1',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novo-naselje/63aec0a83c109c18d20bbe5a',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/nova-detelinara/63aeb4ae568241ab88032ed5',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-oslobodjenja/63ac44313036173803099566',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/centar-novi-sad/63ac7a082af8a275fc0fbaf8',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/gunduliceva/63a45a77d20f9e794102a92a',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novo-naselje/63aedb599b08bc373508ac6a',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/telep/63ad8a1e40c4f3781a04c17e',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/mise-dimitrijevica/6339fa175936bf895b08c5c1',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/socijalno/63ac37f0413a0420b701da8c',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/hadzic-svetic-12/63b3390989347ce492069593',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/dr-djure-jovanovica-kindera/63adbf9eb3fd2529f50c7b0f',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-oslobodjenja/63ac3671413a0420b701d9f2',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/majevicka/63ac739f3240cb1363067cae',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/petra-drapsina/63ab591ab76016982f01339a',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/liman-4/63adb385d4256f04380c1245',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/janka-cmelika/63aebcb1568241ab88032f6d',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novi-sad-centar/63aee98af6dcb2da2701e68e',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/djurdja-brankovica-10/63ae1e88c1c503f48e0f6c41',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/ilije-bircanina/63adeeb282a8172bae02daff',
    'https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/centar-novi-sad/63aec8a12562e3893a016eb4'
)
for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
$ws.Range("A251").ClearContents()
$ws.Range("L243").Select() | Out-Null
